$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Save" header column in H1, reusing the same header style as the
# existing "sum" header in G1 (bold font, thin border, centered alignment).
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Populate the Save column values for the two data rows.
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 1
